$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New cell values in rows 6, 7 (domain-knowledge dataset links/labels) ---
$ws.Range("E6").Value = "https://arxiv.org/pdf/1904.11694.pdf"
$ws.Range("F6").Value = "Domain Knowledge - Data Set"

$ws.Range("E7").Value = "https://arxiv.org/pdf/1906.03523.pdf"
$ws.Range("F7").Value = "Domain Knowledge - Data Set"
$ws.Range("H7").Value = "https://arxiv.org/pdf/2102.11529.pdf"

# --- Row 8: FOL Domain Knowledge label + convert D8/E8 into real hyperlinks ---
$ws.Range("F8").Value = "FOL Domain Knowledge - Data Set"
$ws.Hyperlinks.Add($ws.Range("E8"), "https://arxiv.org/pdf/1805.10872.pdf") | Out-Null

# --- Row 9: Graph data label + hyperlink on E9 ---
$ws.Hyperlinks.Add($ws.Range("E9"), "https://paperswithcode.com/method/mpnn") | Out-Null
$ws.Range("F9").Value = "Graph Data - Domain Knowledge - Data Set"

# D8 hyperlink added last so relationship-id ordering matches rId1..rId3 = E8,E9,D8
$ws.Hyperlinks.Add($ws.Range("D8"), "https://github.com/ML-KULeuven/deepproblog") | Out-Null

# --- Column widths (best-fit approximations for the new layout) ---
$ws.Columns.Item(1).ColumnWidth = 7.421822916666667
$ws.Columns.Item(2).ColumnWidth = 21.257291666666667
$ws.Columns.Item(3).ColumnWidth = 6.586979166666667
$ws.Columns.Item(4).ColumnWidth = 78.92213541666668
$ws.Columns.Item(5).ColumnWidth = 114.58697916666667
$ws.Columns.Item(6).ColumnWidth = 38.592447916666664
$ws.Columns.Item(7).ColumnWidth = 7.257291666666667
$ws.Columns.Item(8).ColumnWidth = 33.42182291666666
$ws.Columns.Item(9).ColumnWidth = 33.42182291666666

# --- Selection moves to E9 ---
$ws.Range("E9").Select() | Out-Null
